# Updated via Streamlit Approval System
#
# The pending-approval sheet used to carry four extra "checkbox" style
# status columns (ACCEPTED / PAID / HOLD / REJECTED) in AP:AS, plus a
# couple of leftover "HOLD" markers in the APPROVAL_1 / APPROVAL_2
# columns (AI/AJ) for the first two data rows. Both are being retired.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the stale "HOLD" markers left in APPROVAL_1 (AI) / APPROVAL_2 (AJ)
# for rows 2 and 3.
$ws.Range("AI2:AJ2").ClearContents()
$ws.Range("AI3:AJ3").ClearContents()

# Remove the now-unused ACCEPTED/PAID/HOLD/REJECTED columns (AP:AS)
# entirely -- this drops the header row cells and every row's boolean
# flag, and shrinks the sheet's used range from A1:AS9 to A1:AO9.
$ws.Columns("AP:AS").Delete()
